$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the first data row (originally 2025-11-15 with no data yet),
# shifting all subsequent rows up by one. This reduces the data range
# from A1:D88 to A1:D87 and makes the former 2025-11-16 row become row 2.
$ws.Rows.Item(2).Delete()
